$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 29164.25
$ws.Range("J3").Value = 29164.25
$ws.Range("L3").Value = 29164.25
$ws.Range("N3").Value = -29392.25
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H48").Value = 4421.4287
$ws.Range("J48").Value = 5000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -15584
$ws.Range("H56").Value = 4421.4287
$ws.Range("J56").Value = 5000
$ws.Range("L56").Value = 15000
$ws.Range("N56").Value = -16068
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488
$ws.Range("H102").Value = 29164.25
$ws.Range("J102").Value = 29164.25
$ws.Range("L102").Value = 29164.25
$ws.Range("N102").Value = -35654.25
$ws.Range("H135").Value = 31251016
$ws.Range("I135").Value = 794.6
$ws.Range("J135").Value = 142858940
$ws.Range("K135").Value = 7151.400000000001
$ws.Range("L135").Value = 1285730460
$ws.Range("M135").Value = -4616.400000000001
$ws.Range("N135").Value = -1285735530
$ws.Range("H136").Value = 40000
$ws.Range("J136").Value = 40000
$ws.Range("L136").Value = 40000
$ws.Range("N136").Value = -50200
$ws.Range("H137").Value = 1589184.1
$ws.Range("I137").Value = 1630.8889
$ws.Range("J137").Value = 3270122.8
$ws.Range("K137").Value = 4892.6667
$ws.Range("L137").Value = 9810368.399999999
$ws.Range("M137").Value = -2342.6667
$ws.Range("N137").Value = -9815468.399999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12823.167
$ws.Range("I32").Value = 7732.256
$ws.Range("J32").Value = 22340.957
$ws.Range("K32").Value = 7732.256
$ws.Range("L32").Value = 22340.957
$ws.Range("M32").Value = -7445.256
$ws.Range("N32").Value = -22914.957
$ws.Range("H61").Value = 1499.0769
$ws.Range("I61").Value = 1198.75
$ws.Range("J61").Value = 1979.6
$ws.Range("K61").Value = 1198.75
$ws.Range("L61").Value = 1979.6
$ws.Range("M61").Value = -986.75
$ws.Range("N61").Value = -2403.6
$ws.Range("H110").Value = 881.5294
$ws.Range("I110").Value = 842.875
$ws.Range("J110").Value = 1500
$ws.Range("K110").Value = 842.875
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = 1202.125
$ws.Range("N110").Value = -5590
$ws.Range("H136").Value = 1499.0769
$ws.Range("I136").Value = 1198.75
$ws.Range("J136").Value = 1979.6
$ws.Range("K136").Value = 3596.25
$ws.Range("L136").Value = 5938.799999999999
$ws.Range("M136").Value = -1046.25
$ws.Range("N136").Value = -11038.8
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 726.6667
$ws.Range("I16").Value = 590
$ws.Range("K16").Value = 590
$ws.Range("M16").Value = -303
$ws.Range("H41").Value = 7291.4
$ws.Range("I41").Value = 3229.5
$ws.Range("J41").Value = 9999.333000000001
$ws.Range("K41").Value = 3229.5
$ws.Range("L41").Value = 9999.333000000001
$ws.Range("M41").Value = -2801.5
$ws.Range("N41").Value = -10855.333
$ws.Range("H113").Value = 726.6667
$ws.Range("I113").Value = 590
$ws.Range("K113").Value = 590
$ws.Range("M113").Value = 1580
$ws.Range("H141").Value = 38317.855
$ws.Range("J141").Value = 38317.855
$ws.Range("L141").Value = 38317.855
$ws.Range("N141").Value = -48677.855
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 749.9143
$ws.Range("I5").Value = 347.15
$ws.Range("K5").Value = 1041.45
$ws.Range("M5").Value = -929.4499999999998
$ws.Range("H92").Value = 481.22726
$ws.Range("I92").Value = 508.5
$ws.Range("J92").Value = 465.64285
$ws.Range("K92").Value = 1525.5
$ws.Range("L92").Value = 1396.92855
$ws.Range("M92").Value = -277.5
$ws.Range("N92").Value = -3892.92855
$ws.Range("H131").Value = 902.16
$ws.Range("I131").Value = 398.75
$ws.Range("J131").Value = 945.93475
$ws.Range("K131").Value = 1196.25
$ws.Range("L131").Value = 2837.80425
$ws.Range("M131").Value = 3843.75
$ws.Range("N131").Value = -12917.80425
$ws.Range("H132").Value = 1192822.6
$ws.Range("I132").Value = 2724.6155
$ws.Range("K132").Value = 24521.5395
$ws.Range("M132").Value = -21991.5395
$ws.Range("H135").Value = 749.9143
$ws.Range("I135").Value = 347.15
$ws.Range("K135").Value = 3124.35
$ws.Range("M135").Value = -589.3499999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 33334668
$ws.Range("I126").Value = 166667260
$ws.Range("J126").Value = 1515.75
$ws.Range("K126").Value = 500001780
$ws.Range("L126").Value = 4547.25
$ws.Range("M126").Value = -499999310
$ws.Range("N126").Value = -9487.25
$ws.Range("H132").Value = 2226814
$ws.Range("I132").Value = 2841.8948
$ws.Range("K132").Value = 8525.6844
$ws.Range("M132").Value = -5995.6844
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 35219.414
$ws.Range("I22").Value = 111550
$ws.Range("J22").Value = 870.65
$ws.Range("K22").Value = 111550
$ws.Range("L22").Value = 870.65
$ws.Range("M22").Value = -111255
$ws.Range("N22").Value = -1460.65
$ws.Range("H27").Value = 35219.414
$ws.Range("I27").Value = 111550
$ws.Range("J27").Value = 870.65
$ws.Range("K27").Value = 111550
$ws.Range("L27").Value = 870.65
$ws.Range("M27").Value = -111443
$ws.Range("N27").Value = -1084.65
$ws.Range("H61").Value = 3170.5881
$ws.Range("I61").Value = 2281.818
$ws.Range("J61").Value = 4800
$ws.Range("K61").Value = 2281.818
$ws.Range("L61").Value = 4800
$ws.Range("M61").Value = -2079.818
$ws.Range("N61").Value = -5204
$ws.Range("H101").Value = 25996.5
$ws.Range("J101").Value = 25996.5
$ws.Range("L101").Value = 25996.5
$ws.Range("N101").Value = -32486.5
$ws.Range("H113").Value = 3170.5881
$ws.Range("I113").Value = 2281.818
$ws.Range("J113").Value = 4800
$ws.Range("K113").Value = 2281.818
$ws.Range("L113").Value = 4800
$ws.Range("M113").Value = -111.8180000000002
$ws.Range("N113").Value = -9140
$ws.Range("H136").Value = 4530.143
$ws.Range("I136").Value = 1326.2
$ws.Range("J136").Value = 12540
$ws.Range("K136").Value = 3978.6
$ws.Range("L136").Value = 37620
$ws.Range("M136").Value = -1428.6
$ws.Range("N136").Value = -42720
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6934
$ws.Range("I136").Value = 2887.652
$ws.Range("J136").Value = 100000
$ws.Range("K136").Value = 8662.956
$ws.Range("L136").Value = 300000
$ws.Range("M136").Value = -6112.956
$ws.Range("N136").Value = -305100
